$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Team N" placeholder labels (merged A3:A7, A8:A12, A13:A14,
# A15:A16) with the actual team-member names. A17:A18 ("team 5") is
# unchanged in content.
$ws.Range("A3").Value = "Shahana, Rahul, Ashish, Amruta"
$ws.Range("A8").Value = "Shreyas"
$ws.Range("A13").Value = "Nandhini"
$ws.Range("A15").Value = "Gaurav"

# Update the view/selection state: scrolled so row 7 is at the top, with
# A8:A12 selected (active cell A8).
$window = $excel.ActiveWindow
$window.ScrollRow = 7
$window.ScrollColumn = 1
$ws.Range("A8:A12").Select()
